# Update time tables for cello and flute camps A and B, and harp camps A and B
# to reflect recent changes in scheduling data: clarify the "Cello MasterClass"
# entries with the teacher's name, correct "Vinnci" -> "Vincci" in the pianist
# rehearsal room name, correct the spelling of Piotr SKWERES's name, and clear
# two stray leftover time values on Day 1-5.

$wb = $excel.ActiveWorkbook

$day1 = $wb.Worksheets.Item("Day 1")
$day2 = $wb.Worksheets.Item("Day 2")
$day3 = $wb.Worksheets.Item("Day 3")
$day4 = $wb.Worksheets.Item("Day 4")
$day5 = $wb.Worksheets.Item("Day 5")

# ---- Cello MasterClass sessions (Day 5) now name the teacher ----
$day5.Range("B7").Value = "Cello MasterClass by Piotr SKWERES`n(Room Piotr)"
$day5.Range("B19").Value = "Cello MasterClass by Piotr SKWERES`n(Room Piotr)"

# ---- Pianist rehearsal room renamed Vinnci -> Vincci (Day 1) ----
$day1.Range("D7").Value = "C3 Rehearsal with pianist`n(Room Vincci)"
$day1.Range("D11").Value = "C1 Rehearsal with pianist`n(Room Vincci) "

# ---- Correct spelling of the teacher's name to "Piotr SWKERES" ----
$day1.Range("C1").Value = "Piotr SWKERES"
$day2.Range("C1").Value = "Piotr SWKERES"
$day3.Range("C1").Value = "Piotr SWKERES"
$day4.Range("C1").Value = "Piotr SWKERES"
$day5.Range("C1").Value = "Piotr SWKERES"

# ---- Clear two stray leftover time values on rows 39-40 (Day 1-5) ----
$day1.Range("A39").ClearContents()
$day1.Range("A40").ClearContents()
$day2.Range("A39").ClearContents()
$day2.Range("A40").ClearContents()
$day3.Range("A39").ClearContents()
$day3.Range("A40").ClearContents()
$day4.Range("A39").ClearContents()
$day4.Range("A40").ClearContents()
$day5.Range("A39").ClearContents()
$day5.Range("A40").ClearContents()

# ---- Reset the stored selection to C1 on each edited sheet ----
$day1.Range("C1").Select()
$day2.Range("C1").Select()
$day3.Range("C1").Select()
$day4.Range("C1").Select()
$day5.Range("C1").Select()
